# "updated js and layout"
# - Agenda slide (slide 2): title placeholder gets an explicit position/size
#   (shrunk height) and the content placeholder gains five new bullet
#   paragraphs ("Question", "Data source", "Data preparation",
#   "Data cleaning", "Tools we used") inserted right after the existing
#   "Introduction and problem statement" line and before "Solution Method".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Title placeholder ("Agenda") - give it an explicit xfrm (same left/top/width
# as inherited from the layout/master, but a reduced height).
$title = $s.Shapes.Item(1)
$title.Left = 66
$title.Top = 28.75
$title.Width = 828
$title.Height = 65.70874015748032

# Content placeholder - rebuild the full bullet list in the new order.
$body = $s.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Introduction and problem statement`rQuestion`rData source`rData preparation`rData cleaning`rTools we used`rSolution Method`rResult"
